$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from serial 45207 to 45208
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = 45208
}
